$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "26.769.37"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "1.595.42"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'210.33"
$ws.Range("E5").Value = "  -2.69%  "
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -1.90%  "
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").Value = "'19.61"
$ws.Range("E10").Value = "  -2.37%  "
$ws.Range("D11").Value = "'0.0836"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").Value = "1.819.16"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("D13").Value = "1.600.68"
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").Value = "'4.05"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "26.742.05"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").Value = "'63.44"
$ws.Range("E17").Value = "  -3.33%  "
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "'209.15"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("D21").Value = "'6.70"
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("D22").Value = "'4.28"
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("E23").Value = "  -7.11%  "
$ws.Range("D24").Value = "'8.85"
$ws.Range("E24").Value = "  -2.88%  "
$ws.Range("D25").Value = "'146.52"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("D26").Value = "'7.47"
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'0.112"
$ws.Range("E28").Value = "  -4.71%  "
$ws.Range("D29").Value = "'15.29"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  -2.65%  "
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("D33").Value = "'0.665"
$ws.Range("E33").Value = "  +22.58%  "
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("D35").Value = "1.311.82"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("D39").Value = "'0.819"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").Value = "'0.788"
$ws.Range("E41").Value = "  -2.23%  "
$ws.Range("E42").Value = "  -3.88%  "
$ws.Range("D43").Value = "'5.29"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "'62.68"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").Value = "1.732.37"
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("D46").Value = "'89.02"
$ws.Range("E46").Value = "  -1.78%  "
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "'0.809"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").Value = "'0.0510"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").Value = "'0.0974"
$ws.Range("E50").Value = "  +2.53%  "
$ws.Range("D51").Value = "'7.46"
$ws.Range("E51").Value = "  -1.15%  "
